$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product names (column B) for rows 2-45
$ws.Range('B2').Value = 'Pink Penguin Puppy Chow'
$ws.Range('B3').Value = 'Corneli Calming Collar™'
$ws.Range('B4').Value = 'YELLOW Homestyle® Adult Dog Food'
$ws.Range('B5').Value = 'Robinson''s Best® Anti Chew Spray'
$ws.Range('B6').Value = 'Doggy Breath Freshener'
$ws.Range('B7').Value = 'Royal® Adult Devplatypus Food'
$ws.Range('B8').Value = 'Anti Fungal Spray'
$ws.Range('B9').Value = 'Caring for your Narwhale by Jenn Petti'
$ws.Range('B10').Value = 'Bird Cage'
$ws.Range('B11').Value = 'Bird Cage Floor Dressing'
$ws.Range('B12').Value = 'Sasquatch Feeder'
$ws.Range('B13').Value = 'Parrot Safety Harness'
$ws.Range('B14').Value = 'Kearney Cat Collar™'
$ws.Range('B15').Value = 'ChowChow® Chewy Treats'
$ws.Range('B16').Value = 'Emo Chow'
$ws.Range('B17').Value = 'AniMarty® Fish Flakes'
$ws.Range('B18').Value = 'Adult Dog Taco Costume'
$ws.Range('B19').Value = 'Barker''s Anti Bark Spray™'
$ws.Range('B20').Value = 'Rose® Collapsible Bowl'
$ws.Range('B21').Value = 'Sasquatch Feed'
$ws.Range('B22').Value = 'Lizard Coat'
$ws.Range('B23').Value = 'Ferret Vitamins'
$ws.Range('B24').Value = 'Car Seat Covers'
$ws.Range('B25').Value = 'Adult Dog T-Rex Costume'
$ws.Range('B26').Value = 'Fox Car Safety Harness'
$ws.Range('B27').Value = 'Dog Agility Tunnel'
$ws.Range('B28').Value = 'Fishpoo® Fish Shampoo'
$ws.Range('B29').Value = 'Carrot Plushie'
$ws.Range('B30').Value = 'Yeti Vitamins'
$ws.Range('B31').Value = 'Freeze-Dried Beef Treats'
$ws.Range('B32').Value = 'Eton Mess® Dog Hoodie'
$ws.Range('B33').Value = 'Rudolph Antler Chews™'
$ws.Range('B34').Value = 'Snake Oil'
$ws.Range('B35').Value = 'Round Fish Bowl'
$ws.Range('B36').Value = 'Parrot Nail Clippers'
$ws.Range('B37').Value = 'Kangaroo® Food Pouch'
$ws.Range('B38').Value = 'Lice & Mite Destroyer'
$ws.Range('B39').Value = 'Fox Toothpaste'
$ws.Range('B40').Value = 'Crunchy Dog Treats'
$ws.Range('B41').Value = 'Mini Adult Dog Food'
$ws.Range('B42').Value = 'All-Purpose Pet Conditioner'
$ws.Range('B43').Value = 'Dental Doggie Treat'
$ws.Range('B44').Value = 'Puppy Snacks'
$ws.Range('B45').Value = 'Dewormer'

# Update Stock and Ordered for row 2
$ws.Range('D2').Value = 5
$ws.Range('E2').Value = 12

$wb.Save()
